# "loading and setting columns by type"
#
# The "data" sheet is a log of objects keyed by ID (col A), with created
# date (col B), modified date (col C, always blank) and name (col D).
# The oldest entry (row 2: "Tzowla Backpack") was removed from the front
# of the log, every remaining row shifted up by one, and a new entry was
# appended at the bottom (new ID, created 12/05/2022, name
# "Tzowla Backpack" again).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..15 up into rows 2..14, one column at a time via Cut so the
# original cell *type* (plain text, incl. the always-empty col C text
# cells) carries over verbatim instead of being re-parsed by Value-assignment
# (which would otherwise "helpfully" reinterpret the date-shaped strings in
# column B as real dates). Column C is intentionally left untouched - every
# cell in it is already an empty text cell, exactly what's needed.
for ($r = 3; $r -le 15; $r++) {
    $dst = $r - 1
    $ws.Range("A" + $r).Cut($ws.Range("A" + $dst))
    $ws.Range("B" + $r).Cut($ws.Range("B" + $dst))
    $ws.Range("D" + $r).Cut($ws.Range("D" + $dst))
}

# New last row: fresh ID and name, column C stays the blank text cell it
# already was.
$ws.Range("A15").Value = "96e75e54-2568-4661-95ca-aa2e68796f19"
$ws.Range("D15").Value = "Tzowla Backpack"

# Writing the new date straight into B15 would get auto-converted to a
# serial date. Stage it as literal text in a scratch cell (forcing text
# format so it isn't parsed as a date), then Cut it into place - Cut moves
# the already-resolved text value/type as-is, and clearing the scratch
# cell's format before removing its row keeps it from leaving any trace.
$scratch = $ws.Range("Z50")
$scratch.NumberFormat = "@"
$scratch.Value = "12/05/2022"
$scratch.Cut($ws.Range("B15"))
$ws.Range("B15").ClearFormats()
$ws.Rows.Item(50).Delete()
